# Applies the cryptos.xlsx update described in the commit diff
# (prices updated, percentage changes updated, and row 31/32 coins swapped)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '40.161.92'
$ws.Cells.Item(2, 5).Value = '  +0.68%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.244.94'
$ws.Cells.Item(3, 5).Value = '  -3.92%  '
$ws.Cells.Item(4, 5).Value = '  -0.08%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '294.15'
$ws.Cells.Item(5, 5).Value = '  -4.92%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '86.31'
$ws.Cells.Item(6, 5).Value = '  +2.75%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.516'
$ws.Cells.Item(7, 5).Value = '  -1.96%  '
$ws.Cells.Item(8, 5).Value = '  -0.05%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.472'
$ws.Cells.Item(9, 5).Value = '  -1.31%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.0800'
$ws.Cells.Item(10, 5).Value = '  -0.34%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '30.67'
$ws.Cells.Item(11, 5).Value = '  +2.90%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '47.88'
$ws.Cells.Item(12, 5).Value = '  -8.71%  '
$ws.Cells.Item(13, 5).Value = '  -2.05%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '6.41'
$ws.Cells.Item(14, 5).Value = '  +0.44%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '2.582.66'
$ws.Cells.Item(15, 5).Value = '  -4.35%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '14.29'
$ws.Cells.Item(16, 5).Value = '  -2.72%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '2.228.41'
$ws.Cells.Item(17, 5).Value = '  -5.52%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.728'
$ws.Cells.Item(18, 5).Value = '  -3.37%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '40.003.40'
$ws.Cells.Item(19, 5).Value = '  +0.24%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '0.0₃0895'
$ws.Cells.Item(20, 5).Value = '  -0.13%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '5.81'
$ws.Cells.Item(21, 5).Value = '  -3.64%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '10.78'
$ws.Cells.Item(22, 5).Value = '  +2.79%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '65.69'
$ws.Cells.Item(23, 5).Value = '  -3.33%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '234.17'
$ws.Cells.Item(24, 5).Value = '  -0.28%  '
$ws.Cells.Item(25, 5).Value = '  +0.09%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '2.43'
$ws.Cells.Item(26, 5).Value = '  -3.83%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '1.86'
$ws.Cells.Item(27, 5).Value = '  +2.71%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '23.18'
$ws.Cells.Item(28, 5).Value = '  -1.10%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.20'
$ws.Cells.Item(29, 5).Value = '  +3.63%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '9.28'
$ws.Cells.Item(30, 5).Value = '  +0.73%  '
$ws.Cells.Item(31, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '33.55'
$ws.Cells.Item(31, 5).Value = '  -2.11%  '
$ws.Cells.Item(32, 2).Value = 'Monero'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '155.33'
$ws.Cells.Item(32, 5).Value = '  +1.71%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.998'
$ws.Cells.Item(33, 5).Value = '  -0.38%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '4.87'
$ws.Cells.Item(34, 5).Value = '  -3.76%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.0714'
$ws.Cells.Item(35, 5).Value = '  +0.13%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.37'
$ws.Cells.Item(36, 5).Value = '  -4.49%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '16.71'
$ws.Cells.Item(37, 5).Value = '  +7.95%  '
$ws.Cells.Item(38, 5).Value = '  -1.02%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.0988'
$ws.Cells.Item(39, 5).Value = '  +0.69%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '2.70'
$ws.Cells.Item(40, 5).Value = '  -1.83%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.68'
$ws.Cells.Item(41, 5).Value = '  -1.24%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '3.79'
$ws.Cells.Item(42, 5).Value = '  +0.80%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '1.959.80'
$ws.Cells.Item(43, 5).Value = '  -0.48%  '
$ws.Cells.Item(44, 5).Value = '  -3.03%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.0271'
$ws.Cells.Item(45, 5).Value = '  +3.15%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '9.54'
$ws.Cells.Item(46, 5).Value = '  +1.61%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '16.44'
$ws.Cells.Item(47, 5).Value = '  -4.97%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '2.63'
$ws.Cells.Item(48, 5).Value = '  -1.02%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '2.448.87'
$ws.Cells.Item(49, 5).Value = '  -4.44%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '71.27'
$ws.Cells.Item(50, 5).Value = '  +2.08%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '1.47'
$ws.Cells.Item(51, 5).Value = '  +8.20%  '
